$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# STEP 8 + STEP 9 values (column E, rows 3-9 and 11-17)
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 288
$ws.Range("E5").Value = 72
$ws.Range("E6").Value = 216
$ws.Range("E7").Value = 72
$ws.Range("E8").Value = 216
$ws.Range("E9").Value = 144

$ws.Range("E11").Value = 0
$ws.Range("E12").Value = 34560
$ws.Range("E13").Value = 8640
$ws.Range("E14").Value = 25920
$ws.Range("E15").Value = 8640
$ws.Range("E16").Value = 25920
$ws.Range("E17").Value = 17280

# E7 no longer carries the custom number-format style (s="1") -> back to default/Normal style
$ws.Range("E7").Style = "Normal"

# Selection moved from A1:E1 to the whole column F (F1:F1048576), active cell F1
$ws.Columns.Item(6).Select()
